$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells for columns I (I0) and J (IF), matching the formatting
# of the existing header cells (bold, centered, bordered) by copying the
# format from H1.
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Data values for the new I and J columns (rows 2-20). In every row the
# two columns hold the same value.
$values = @{
    2  = 8
    3  = 8
    4  = 7
    5  = 8
    6  = 6
    7  = 6
    8  = 7
    9  = 6
    10 = 6
    11 = 7
    12 = 7
    13 = 7
    14 = 8
    15 = 4
    16 = 7
    17 = 5
    18 = 4
    19 = 3
    20 = 7
}

foreach ($row in $values.Keys) {
    $val = $values[$row]
    $ws.Cells.Item($row, 9).Value = $val
    $ws.Cells.Item($row, 10).Value = $val
}
